$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-content permutation for rows 114-127 (species-record shuffle), derived from the
# target diff: the record identity/location columns (A,B,D,E,F,G,H,M,Q,R,AJ,AK,AO) move
# between fixed row numbers while all other metadata columns stay put.

$ws.Range("A114").Value = 111743554
$ws.Range("B114").Value = 88966
$ws.Range("D114").Value = "NT"
$ws.Range("E114").Value = 5754
$ws.Range("F114").Value = "Gultoppig fingersvamp"
$ws.Range("G114").Value = "Ramaria testaceoflava"
$ws.Range("H114").Value = "(Bres.) Corner"
$ws.Range("Q114").Value = 339577.2032005055
$ws.Range("R114").Value = 6571127.007499221
$ws.Range("A115").Value = 111743520
$ws.Range("B115").Value = 56398
$ws.Range("D115").Value = "NT"
$ws.Range("E115").Value = 100109
$ws.Range("F115").Value = "Tretåig hackspett"
$ws.Range("G115").Value = "Picoides tridactylus"
$ws.Range("H115").Value = "(Linnaeus, 1758)"
$ws.Range("M115").Value = "färska spår"
$ws.Range("Q115").Value = 339096.8530521042
$ws.Range("R115").Value = 6571013.66294401
$ws.Range("AJ115").Value = "gran"
$ws.Range("AK115").Value = "Picea abies"
$ws.Range("AO115").Value = "Picea abies"
$ws.Range("A116").Value = 111743523
$ws.Range("B116").Value = 73634
$ws.Range("D116").Value = "LC"
$ws.Range("E116").Value = 6426
$ws.Range("F116").Value = "Kattfotslav"
$ws.Range("G116").Value = "Felipes leucopellaeus"
$ws.Range("H116").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q116").Value = 339009.0243061834
$ws.Range("R116").Value = 6571011.238422027
$ws.Range("A118").Value = 111743524
$ws.Range("B118").Value = 94134
$ws.Range("D118").Value = "NT"
$ws.Range("E118").Value = 53
$ws.Range("F118").Value = "Vedtrappmossa"
$ws.Range("G118").Value = "Crossocalyx hellerianus"
$ws.Range("H118").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q118").Value = 338949.7235384365
$ws.Range("R118").Value = 6571040.381812023
$ws.Range("AJ118").Value = "tall"
$ws.Range("AK118").Value = "Pinus sylvestris"
$ws.Range("AO118").Value = "Pinus sylvestris"
$ws.Range("A119").Value = 111743515
$ws.Range("Q119").Value = 339441.7613444271
$ws.Range("R119").Value = 6571017.506567059
$ws.Range("A120").Value = 111743517
$ws.Range("B120").Value = 73634
$ws.Range("D120").Value = "LC"
$ws.Range("E120").Value = 6426
$ws.Range("F120").Value = "Kattfotslav"
$ws.Range("G120").Value = "Felipes leucopellaeus"
$ws.Range("H120").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q120").Value = 339278.3213300391
$ws.Range("R120").Value = 6571107.378548244
$ws.Range("A121").Value = 111743516
$ws.Range("B121").Value = 96348
$ws.Range("D121").Value = "VU"
$ws.Range("E121").Value = 220787
$ws.Range("F121").Value = "Knärot"
$ws.Range("G121").Value = "Goodyera repens"
$ws.Range("H121").Value = "(L.) R. Br."
$ws.Range("Q121").Value = 339415.5147437509
$ws.Range("R121").Value = 6571015.54325202
$ws.Range("A122").Value = 111743519
$ws.Range("B122").Value = 90666
$ws.Range("D122").Value = "LC"
$ws.Range("E122").Value = 4364
$ws.Range("F122").Value = "Dropptaggsvamp"
$ws.Range("G122").Value = "Hydnellum ferrugineum"
$ws.Range("H122").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q122").Value = 339118.4126724883
$ws.Range("R122").Value = 6571062.424656671
$ws.Range("A123").Value = 111743551
$ws.Range("B123").Value = 96348
$ws.Range("D123").Value = "VU"
$ws.Range("E123").Value = 220787
$ws.Range("F123").Value = "Knärot"
$ws.Range("G123").Value = "Goodyera repens"
$ws.Range("H123").Value = "(L.) R. Br."
$ws.Range("Q123").Value = 339522.8608171764
$ws.Range("R123").Value = 6571091.407599592
$ws.Range("A124").Value = 111743526
$ws.Range("B124").Value = 90666
$ws.Range("D124").Value = "LC"
$ws.Range("E124").Value = 4364
$ws.Range("F124").Value = "Dropptaggsvamp"
$ws.Range("G124").Value = "Hydnellum ferrugineum"
$ws.Range("H124").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q124").Value = 338870.1217119552
$ws.Range("R124").Value = 6571086.774471543
$ws.Range("A125").Value = 111743546
$ws.Range("Q125").Value = 339474.5644867857
$ws.Range("R125").Value = 6571113.931964876
$ws.Range("A127").Value = 111743527
$ws.Range("B127").Value = 96348
$ws.Range("D127").Value = "VU"
$ws.Range("E127").Value = 220787
$ws.Range("F127").Value = "Knärot"
$ws.Range("G127").Value = "Goodyera repens"
$ws.Range("H127").Value = "(L.) R. Br."
$ws.Range("Q127").Value = 338598.1684531783
$ws.Range("R127").Value = 6571109.585305012

# Cells that no longer have content after the shuffle (the donor row for these target rows
# did not have a value in these optional columns).
$ws.Range("AJ124").ClearContents()
$ws.Range("AK124").ClearContents()
$ws.Range("AO124").ClearContents()
$ws.Range("M127").ClearContents()
$ws.Range("AJ127").ClearContents()
$ws.Range("AK127").ClearContents()
$ws.Range("AO127").ClearContents()
